# ---------------------------------------------------------------------------
# 20240525 - Data Science Personal Log.xlsx
# "added hockey and updated movies"
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# ===========================================================================
# 1. Movies sheet - fix up row 7 and append three new rows (8-10)
# ===========================================================================
$movies = $wb.Worksheets.Item("Movies")

# Row 7 edits: Godzilla Minus One rating/genre/comment update
$movies.Range("D7").Value = 7.4
$movies.Range("E7").Value = "Drama/Action"
$movies.Range("F7").Value = "pretty good. Impressive work and creative story. You think the japanese could have mustered more than one plane though."

# New row 8
$movies.Range("B8").Value = 45454
$movies.Range("B8").NumberFormat = "d-mmm-yy"
$movies.Range("C8").Value = "Godzilla.X.Kong.The.New.Empire.2024"
$movies.Range("D8").Value = 5.3
$movies.Range("E8").Value = "Action"
$movies.Range("F8").Value = "bad. Stupid. Did not finish. I was out when there was a scottish man driving a spaceship into some wormhole thing"

# New row 9
$movies.Range("B9").Value = 45455
$movies.Range("B9").NumberFormat = "d-mmm-yy"
$movies.Range("C9").Value = "Ravenous"
$movies.Range("D9").Value = 9.2
$movies.Range("E9").Value = "Horror"
$movies.Range("F9").Value = "guy pearce and danny bowel was awesome. Cannabalism stuff was really sick"

# New row 10
$movies.Range("B10").Value = 45455
$movies.Range("B10").NumberFormat = "d-mmm-yy"
$movies.Range("C10").Value = "Postcard Killings"
$movies.Range("D10").Value = 5.4
$movies.Range("E10").Value = "Drama/Triller"
$movies.Range("F10").Value = "not very good. Poor villans, not very belivable or scary"

# ===========================================================================
# 2. Videogames sheet - brand new, inserted right after "Shows"
# ===========================================================================
$shows = $wb.Worksheets.Item("Shows")
$videogames = $wb.Worksheets.Add($null, $shows)
$videogames.Name = "Videogames"

$videogames.Range("C4").Value = "Name"
$videogames.Range("D4").Value = "Year of Release"
$videogames.Range("E4").Value = "Year of Completion"
$videogames.Range("F4").Value = "Difficulty"
$videogames.Range("G4").Value = "Enjoyment"
$videogames.Range("H4").Value = "Badass?"
$videogames.Range("I4").Value = "Finished?"
$videogames.Range("J4").Value = "Keep or Pawn?"
$videogames.Range("K4").Value = "Comment"

$videogames.Range("C5").Value = "Robocop"
$videogames.Range("D5").Value = 2023
$videogames.Range("E5").Value = 2024
$videogames.Range("F5").Value = 6.5
$videogames.Range("G5").Value = 8.9
$videogames.Range("H5").Value = "Y"
$videogames.Range("I5").Value = "Y"
$videogames.Range("J5").Value = "Keep"
$videogames.Range("K5").Value = "Surpringly great game. Very fun"

$videogames.Range("G32").Select() | Out-Null

# ===========================================================================
# 3. Poutine sheet - cosmetic column-width tweak on column F
# ===========================================================================
$poutine = $wb.Worksheets.Item("Poutine")
$poutine.Columns.Item(6).ColumnWidth = 6.6

# ===========================================================================
# 4. Books sheet - add the header row
# ===========================================================================
$books = $wb.Worksheets.Item("Books")
$books.Range("C5").Value = "Date"
$books.Range("D5").Value = "Name"
$books.Range("E5").Value = "Genre"
$books.Range("F5").Value = "Audiobook?"
$books.Range("G5").Value = "Rating"
$books.Range("H5").Value = "Comment"
$books.Range("C6").Select() | Out-Null

# ===========================================================================
# 5. Hockey sheet - fill in row 6 with the McCormick game details
# ===========================================================================
$hockey = $wb.Worksheets.Item("Hockey")
$hockey.Range("C6").Value = "McCormick"
$hockey.Range("E6").Value = 1
$hockey.Range("F6").Value = 2
$hockey.Range("G6").Value = "wrist shot"
$hockey.Range("H6").Value = "felt slow"
$hockey.Range("I6").Value = "passing needs work. Stop eating so much before games. Tired because did some biking and moed both lawns"

# ===========================================================================
# 6. Working Out sheet - brand new, appended after "Hockey"
# ===========================================================================
$workingout = $wb.Worksheets.Add($null, $hockey)
$workingout.Name = "Working Out"

$workingout.Range("C3").Value = "Date"
$workingout.Range("D3").Value = "Plank"
$workingout.Range("E3").Value = "Pushup"
$workingout.Range("F3").Value = "Headstand"
$workingout.Range("G3").Value = "Boxing"
$workingout.Range("H3").Value = "Yoga"

$workingout.Range("C4").Value = 41084
$poutine.Range("B5").Copy() | Out-Null
$workingout.Range("C4").PasteSpecial(-4122) | Out-Null
$workingout.Range("F4").Value = "y"

$workingout.Columns.Item(3).ColumnWidth = 9.73046875

$workingout.Range("J12").Select() | Out-Null

# ===========================================================================
# 7. Final view state - Hockey tab active, selections on each sheet
# ===========================================================================
$movies.Activate()
$movies.Range("I20").Select() | Out-Null

$shows.Activate()
$shows.Range("H21").Select() | Out-Null

$poutine.Activate()
$poutine.Range("J36").Select() | Out-Null

$hockey.Activate()
$hockey.Range("I7").Select() | Out-Null
